$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "آموزش رایگان اقتصاد مهندسی | مکتب‌خونه"
$ws.Range("B2").Value = "نوید خادمی"
$ws.Range("C2").Value = "دانشگاه تهران"
$ws.Range("D2").Value = "رایگان"
$ws.Range("E2").Value = "14 جلسه"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://maktabkhooneh.orghttps://maktabkhooneh.org/course/%D8%A2%D9%85%D9%88%D8%B2%D8%B4-%D8%B1%D8%A7%DB%8C%DA%AF%D8%A7%D9%86-%D8%A7%D9%82%D8%AA%D8%B5%D8%A7%D8%AF-%D9%85%D9%87%D9%86%D8%AF%D8%B3%DB%8C-mk1364/")

$ws.Range("A3").Value = "آموزش مدیریت سرور و امنیت در لینوکس | مکتب‌خونه"
$ws.Range("B3").Value = "Greg Williams"
$ws.Range("C3").Value = "مکتب‌خونه"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49,000"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "3 ساعت"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://maktabkhooneh.orghttps://maktabkhooneh.org/course/%D8%A2%D9%85%D9%88%D8%B2%D8%B4-%D9%85%D8%AF%DB%8C%D8%B1%DB%8C%D8%AA-%D8%B3%D8%B1%D9%88%D8%B1-%D8%A7%D9%85%D9%86%DB%8C%D8%AA-%D9%84%DB%8C%D9%86%D9%88%DA%A9%D8%B3-mk1330/")
